# Incorporación del distrito federal 2015 en xlsx nombre_elecciones
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 3 (after the "Presidencia Municipal 15" row),
# shifting everything else down, and fill it with the new "Distrito federal 15" entry.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "Distrito federal 15"
$ws.Range("B3").Value = "df_15"
$ws.Range("C3").Value = "#ffc8dd"

$ws.Range("F8").Select()
